# feat: add 2022-Q4 data
#
# The workbook tracks quarterly holdings. A new "2022-Q4" sheet of data is
# introduced, which pushes the existing quarters down:
#   - the data that used to live on the "2022-Q3" sheet becomes the new
#     "2022-Q4" sheet content (refreshed with the latest figures),
#   - a fresh copy of the former "2022-Q3" sheet (with its original,
#     untouched figures) is inserted right after it and becomes the new
#     historical "2022-Q3" sheet,
#   - the "2022-Q2" sheet is left as-is (it simply shifts one position to
#     the right),
#   - the "总计" (summary) sheet gains a new row for 2022-Q2 and its
#     existing rows are relabeled to 2022-Q4 / 2022-Q3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the current "2022-Q3" sheet. The duplicate preserves the
#    original Q3 numbers and will be renamed to stay "2022-Q3", while the
#    original sheet object is repurposed (renamed + refreshed) to become
#    "2022-Q4".
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Copy($null, $q3Sheet)

$q3Copy = $wb.Worksheets.Item(3)

# Free up the "2022-Q3" name before assigning it to the copy.
$q3Sheet.Name = "2022-Q4"
$q3Copy.Name = "2022-Q3"

# Helper: assign a value that must be stored as text (these columns hold
# numeric-looking strings, not real numbers) without leaving a lingering
# "quote prefix" style on the cell.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 2. Refresh the figures on the new "2022-Q4" sheet.
# ---------------------------------------------------------------------
Set-TextValue $q3Sheet.Range("D2") "0.21"
Set-TextValue $q3Sheet.Range("E2") "86.58"
Set-TextValue $q3Sheet.Range("F2") "5.84"
Set-TextValue $q3Sheet.Range("G2") "0.0123"

Set-TextValue $q3Sheet.Range("D3") "0.16"
Set-TextValue $q3Sheet.Range("E3") "86.58"
Set-TextValue $q3Sheet.Range("F3") "5.84"
Set-TextValue $q3Sheet.Range("G3") "0.0093"

# Keep the "2022-Q2" sheet as the selected/active tab, matching its
# original state.
$wb.Worksheets.Item(4).Activate()

# ---------------------------------------------------------------------
# 3. Update the "总计" summary sheet: relabel the existing two rows and
#    append a new row for 2022-Q2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("B3").Value = "2022-Q3"

# Copy formatting (bold font + border) from row 3's A cell into row 4's A
# cell before writing its value.
$summary.Range("A3").Copy($summary.Range("A4"))
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.02

Write-Host "Applied 2022-Q4 rollover edits"
